# fix: Force Mastercard Orange headings
# Sets heading run colors to #FF5F00 (Mastercard Orange) for every
# paragraph using a "Heading*" style in the document body, and updates
# the underlying Heading5/Heading6/Title/Subtitle paragraph styles
# (in styles.xml) so the theme-color overrides are replaced with the
# explicit RGB value as well.

$d = $word.ActiveDocument

# Mastercard Orange FF5F00, expressed as a Word COM BGR-packed color
# value (R | G<<8 | B<<16): 0xFF | 0x5F<<8 | 0x00<<16 = 24575.
$orange = 24575

# 1) Color the run text of every Heading1/Heading2/Heading3 paragraph in
#    the document body. Only the run's rPr gets the color (not the
#    paragraph mark), so build a Range that stops one character short of
#    the paragraph end (which excludes the trailing paragraph mark).
foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -like "Heading*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1
        if ($end -gt $start) {
            $r = $d.Range($start, $end)
            $r.Font.Color = $orange
        }
    }
}

# 2) Force the heading-adjacent styles that still carried theme-color
#    overrides (Heading5, Heading6, Title, Subtitle) to the same
#    explicit orange, replacing their theme color/shade references.
foreach ($styleName in @("Heading5", "Heading6", "Title", "Subtitle")) {
    $s = $d.Styles($styleName)
    $s.Font.Color = $orange
}
